$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Mark "Read Test Passed" column (C) as TRUE for all data rows (2-24)
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 3).Value = $true
}

# Update the view state: scroll so row 11 is at the top and select C21
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("C21").Select()
